$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.857.47"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.638.93"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").Value = "216.96"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "19.86"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.866.99"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "1.643.30"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "67.21"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "26.842.08"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "218.16"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "2.44"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("D24").Value = "9.15"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").Value = "147.15"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "15.77"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").Value = "1.265.20"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "0.835"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "1.778.12"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "61.95"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "91.76"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "1.59"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").Value = "  -0.81%  "
